$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20: add hours worked (D20) and task description (E20)
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = "mulighet for sletting av comments for author og admin"

# Row 21: add hours worked (D21) and task description (E21)
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = "mulighet for edit av content details for owner og admin"

# Column E widened (best fit) to accommodate the new, longer text
$ws.Columns.Item(5).ColumnWidth = 49.67

# Move active cell selection to reflect where the author was last working
$ws.Range("G16").Select()
